$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.763.32"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.797.04"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.28"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4563"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3715"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07241"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8560"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.41"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.797.57"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.308"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.508"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07032"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.29"
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9991"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.62"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.756.22"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.290"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.019.72"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.51"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.161"
$ws.Range("E27").Value = "  -12.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.204"
$ws.Range("E29").Value = "  -3.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.14"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08838"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7587"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.446"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.885"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9985"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.113"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01941"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05220"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.388"
$ws.Range("E40").Value = "  +4.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.894"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.132"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5239"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1647"
$ws.Range("E44").Value = "  -5.12%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5040"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.489"
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.29"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9987"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.648"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06284"
$ws.Range("E51").Value = "  -1.58%  "
